$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'273.31"
$ws.Range("D3").Value = "'26.84"
$ws.Range("E3").Value = "'0.41%"
$ws.Range("D4").Value = "'4.904"
$ws.Range("E4").Value = "'3.87%"
$ws.Range("D5").Value = "'0.06320"
$ws.Range("E5").Value = "'3.14%"
$ws.Range("D6").Value = "'6.908"
$ws.Range("E6").Value = "'2.35%"
$ws.Range("D7").Value = "'3.355"
$ws.Range("E7").Value = "'5.44%"
$ws.Range("D8").Value = "'1.364"
$ws.Range("E8").Value = "'51.73%"
$ws.Range("D9").Value = "'0.8831"
$ws.Range("E9").Value = "'3.24%"
$ws.Range("E10").Value = "'2.32%"
$ws.Range("D11").Value = "'0.05092"
$ws.Range("E11").Value = "'0.07%"
$ws.Range("D12").Value = "'0.07407"
$ws.Range("E12").Value = "'3.50%"
$ws.Range("D13").Value = "'0.03148"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("D14").Value = "'0.09039"
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("D15").Value = "'0.001565"
$ws.Range("E15").Value = "'2.38%"
$ws.Range("D16").Value = "'0.0006314"
$ws.Range("E16").Value = "'3.63%"
$ws.Range("D17").Value = "'0.006006"
$ws.Range("E17").Value = "'-1.18%"
$ws.Range("D18").Value = "'3.467"
$ws.Range("E18").Value = "'0.07%"
$ws.Range("E19").Value = "'-0.26%"
$ws.Range("E20").Value = "'2.50%"
$ws.Range("D21").Value = "'0.1332"
$ws.Range("E21").Value = "'3.95%"
$ws.Range("D22").Value = "'3.912"
$ws.Range("E22").Value = "'1.83%"
$ws.Range("D23").Value = "'0.04338"
$ws.Range("E23").Value = "'2.30%"
$ws.Range("D24").Value = "'0.001177"
$ws.Range("E24").Value = "'0.02%"
$ws.Range("D25").Value = "'0.003652"
$ws.Range("E25").Value = "'-12.04%"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.12%"
$ws.Range("E27").Value = "'1.14%"
$ws.Range("D40").Value = "'0.04044"
$ws.Range("E40").Value = "'1.88%"
$ws.Range("D41").Value = "'0.006624"
$ws.Range("E41").Value = "'57.76%"
$ws.Range("D42").Value = "'0.1164"
$ws.Range("E42").Value = "'3.89%"
$ws.Range("D43").Value = "'0.002130"
$ws.Range("E43").Value = "'4.44%"
$ws.Range("D44").Value = "'0.01255"
$ws.Range("E44").Value = "'7.30%"
$ws.Range("D45").Value = "'0.00005333"
$ws.Range("E45").Value = "'3.35%"
$ws.Range("E46").Value = "'159.93%"
$ws.Range("E47").Value = "'-29.20%"
